# "Finish Market Random Call"
# - Delete Sheet3 (unused empty sheet)
# - market_item sheet: add a "tag" column (G) with header + "market" value on
#   every data row, and add 5 more market items (normalAttack3..normalAttack7)
# - market_level_resource sheet content is unchanged, only selection moves

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Remove Sheet3 (empty placeholder sheet)
# ---------------------------------------------------------------------------
$wb.Worksheets("Sheet3").Delete() | Out-Null

# ---------------------------------------------------------------------------
# 2. market_item: new "tag" column header + fill existing rows, append rows
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets("market_item")

# New header cells for column G
$ws1.Cells.Item(1, 7).Value = "market 物品商店 ability 技能商店"
$ws1.Cells.Item(2, 7).Value = "tag"

# Fill "tag" for the two pre-existing data rows
$ws1.Cells.Item(3, 7).Value = "market"
$ws1.Cells.Item(4, 7).Value = "market"

# New data rows 5-9 (normalAttack3..normalAttack7), each tagged "market"
$newItems = @(
    @(3, "normalAttack3", "market_item_CN_normalAttack3"),
    @(4, "normalAttack4", "market_item_CN_normalAttack4"),
    @(5, "normalAttack5", "market_item_CN_normalAttack5"),
    @(6, "normalAttack6", "market_item_CN_normalAttack6"),
    @(7, "normalAttack7", "market_item_CN_normalAttack7")
)

$r = 5
foreach ($item in $newItems) {
    $ws1.Cells.Item($r, 1).Value = $item[0]
    $ws1.Cells.Item($r, 2).Value = $item[1]
    $ws1.Cells.Item($r, 3).Value = $item[2]
    $ws1.Cells.Item($r, 4).Value = "test_gloves.png"
    $ws1.Cells.Item($r, 5).Value = 1
    $ws1.Cells.Item($r, 6).Value = 1
    $ws1.Cells.Item($r, 7).Value = "market"
    $r++
}

# Widen the new tag column and restore the outline-level bookkeeping Excel
# keeps on the sheet (mirrors the author's interactive column grouping).
$ws1.Columns.Item(7).Group() | Out-Null
$ws1.Columns.Item(7).Group() | Out-Null
$ws1.Columns.Item(7).Group() | Out-Null
$ws1.Columns.Item(7).Group() | Out-Null
$ws1.Columns.Item(7).Group() | Out-Null
$ws1.Columns.Item(7).Group() | Out-Null
$ws1.Columns.Item(7).ColumnWidth = 54.15

# ---------------------------------------------------------------------------
# 3. market_level_resource: no value changes, just move the active cell
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets("market_level_resource")
$ws2.Range("C1").Select() | Out-Null

# ---------------------------------------------------------------------------
# 4. Re-activate market_item (it stays the selected tab) and restore its
#    final active cell.
# ---------------------------------------------------------------------------
$ws1.Activate() | Out-Null
$ws1.Range("G16").Select() | Out-Null
